$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.959.67"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "3.078.15"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "3.074.45"
$ws.Range("E9").Value = "  -1.60%  "
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000242"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.52%  "
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").Value = "3.588.70"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").Value = "66.852.75"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.36%  "
$ws.Range("D20").Value = "3.075.12"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "486.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.690"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.67%  "
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.112"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.29%  "
$ws.Range("D34").Value = "0.0₃0910"
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.67"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.952"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "46.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("E40").Value = "  -4.44%  "
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.31%  "
$ws.Range("D43").Value = "2.762.94"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "370.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "136.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0344"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("E51").Value = "  -1.68%  "
